# Apply "add product page complete" edits to the Products sheet.
# Updates Part No (P), Brand (C), Model (D), Type (G), Range (H) and
# Description (N) values per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# --- Row 2 : PRESSURE SWITCH (H-302) ---
$ws.Range("P2").Value2 = "PS-0.15 BAR"

# --- Row 3 : BMS / ETAMATIC ---
$ws.Range("P3").Value2 = "BMS-ETAMATIC-OEM-110V"

# --- Row 4 : FLAME SCANNER ---
$ws.Range("P4").Value2 = "FLAME-SNSR-FF207"

# --- Row 5 : PILOT IGNITION ---
$ws.Range("P5").Value2 = "IG-ROD-PILOT"

# --- Row 6 : IGNITION TRANSFORMER ---
$ws.Range("P6").Value2 = "IG-XF-110V"

# --- Row 7 : PRESSURE SWITCH (TRAFAG, min oil) ---
$ws.Range("P7").Value2 = "PS-MIN-OIL"

# --- Row 8 : PRESSURE SWITCH (TRAFAG, max oil) ---
$ws.Range("P8").Value2 = "PS-MAX-OIL"

# --- Row 9 : IO CARD - digital input ---
$ws.Range("P9").Value2 = "PLC-DI-1734-IB8"

# --- Row 10 : IO CARD - analog output ---
$ws.Range("P10").Value2 = "PLC-AO-1734-OE4C"

# --- Row 11 : IO CARD - digital output ---
$ws.Range("P11").Value2 = "PLC-DO-1734-OB8"

# --- Row 12 : IO CARD - analog input ---
$ws.Range("P12").Value2 = "PLC-AI-1734-IE8C"

# --- Row 13 : IO CARD BACK-PANE ---
$ws.Range("N13").Value2 = "H-302 HEATER PANEL SPARE"
$ws.Range("P13").Value2 = "PLC-BACKPANE-1734-TOP"

# --- Row 14 : SOLENOID COIL ---
$ws.Range("P14").Value2 = "SOV-DIESEL"

# --- Row 15 : 8 PIN RELAY (SIEMENS) ---
$ws.Range("P15").Value2 = "RELAY-SIEMENS-7RQ0201-1BX00-240VAC"

# --- Row 16 : 8 PIN RELAY (OMRON) ---
$ws.Range("P16").Value2 = "RELAY-OMRON-220VAC"

# --- Row 17 : SERVO MOTOR (ROLOFF) ---
$ws.Range("N17").Value2 = "H-302 HEATER SPARE"
$ws.Range("P17").Value2 = "RELAY-ROLOFF-115VAC"

# --- Row 18 : BURNER CONTROLLER (HONEYWELL) ---
$ws.Range("N18").Value2 = "HONEYWELL WITH AMPLIFIER CARD R7847"
$ws.Range("P18").Value2 = "BMS-HONEYWELL-EC7823"

# --- Row 19 : IGNITION TRANSFORMER (DONGAN) ---
$ws.Range("P19").Value2 = "IG-XF-DONGAN-220VAC"

# --- Row 20 : BURNER CONTROLLER (RATIOTRONIC 6006) ---
$ws.Range("C20").Value2 = "RATIOTRONIC"
# "6006" looks numeric, force the cell to remain Text like the source file
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "6006"
$ws.Range("P20").Value2 = "BLR-CONTROLLER-RATIOTRONIC-6006-110 VAC/ 250 VAC"

# --- Row 21 : SERVO MOTOR (RATIOTRONIC 6026) ---
$ws.Range("C21").Value2 = "RATIOTRONIC"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "6026"
$ws.Range("P21").Value2 = "BLR-ACTUATOR-RATIOTRONIC-6026---"

# --- Row 22 : PRESSURE SWITCH (RATIOTRONIC) ---
$ws.Range("C22").Value2 = "RATIOTRONIC"
$ws.Range("G22").Value2 = "PS"
$ws.Range("P22").Value2 = "BLR-PS-RATIOTRONIC-GGW 150 A4/2X-500 mbar"

# --- Row 23 : LEVEL CONTROLLER (RATIOTRONIC) ---
$ws.Range("C23").Value2 = "RATIOTRONIC"
$ws.Range("P23").Value2 = "BLR-CONTROLLER-RATIO TRONIC-LC3050-220VAC"

# --- Row 24 : VFD DISPLAY ---
$ws.Range("H24").Value2 = ""
$ws.Range("N24").Value2 = "BOILER SPARE ITEM"
$ws.Range("P24").Value2 = "BLR-DISPLAY-DUNPHY-ETCRAT00755-"

# --- Row 25 : PRESSURE TRANSMITTER ---
$ws.Range("H25").Value2 = "10 bar"
$ws.Range("N25").Value2 = "BOILER SPARE ITEM"
$ws.Range("P25").Value2 = "BLR-PT-WIKA-A10-10 bar"

# --- Row 26 : SOLENOID VALVE ---
$ws.Range("P26").Value2 = "BLR-SOLENOID-SIEMENS-321H2322-30 BAR"

# --- Row 27 : OVERLOAD RELAY ---
$ws.Range("N27").Value2 = "BOILER SPARE PART"
$ws.Range("P27").Value2 = "OLR-PARKER-10A"

# --- Row 28 : MAGNETIC CONTACTOR ---
$ws.Range("P28").Value2 = "CONTACTOR-SIEMENS-3RT2016-1APO1"

# --- Row 29 : IGNITION TRANSFORMER (SIEMENS) ---
$ws.Range("P29").Value2 = "IG-XF-SIEMENS-ZA 20 100 LH 21"

# --- Row 30 : IGNITION ROD ---
$ws.Range("P30").Value2 = "BLR-IG-ROD"
